# Auto-generated edit script: updates cryptos list Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.516.97'
$ws.Range("E2").Value = '  +1.63%  '
$ws.Range("D3").Value = '3.942.94'
$ws.Range("E3").Value = '  +0.25%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '506.44'
$ws.Range("E5").Value = '  +4.11%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.93'
$ws.Range("E6").Value = '  +0.04%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.626'
$ws.Range("E7").Value = '  -0.18%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.736'
$ws.Range("E9").Value = '  +0.17%  '
$ws.Range("E10").Value = '  +4.00%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0000351'
$ws.Range("E11").Value = '  -1.07%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '43.64'
$ws.Range("E12").Value = '  +1.46%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.52'
$ws.Range("E13").Value = '  -1.71%  '
$ws.Range("D14").Value = '4.571.48'
$ws.Range("D15").Value = '3.953.71'
$ws.Range("E15").Value = '  +0.58%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.26'
$ws.Range("E16").Value = '  -1.92%  '
$ws.Range("E17").Value = '  -0.22%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.22'
$ws.Range("E18").Value = '  +7.21%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '20.04'
$ws.Range("E19").Value = '  +0.22%  '
$ws.Range("D20").Value = '69.606.83'
$ws.Range("E20").Value = '  +1.60%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '436.74'
$ws.Range("E21").Value = '  -1.22%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.45'
$ws.Range("E22").Value = '  -1.77%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.71'
$ws.Range("E23").Value = '  -2.65%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '89.24'
$ws.Range("E24").Value = '  +0.90%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.98'
$ws.Range("E25").Value = '  +6.23%  '
$ws.Range("E26").Value = '  +6.80%  '
$ws.Range("E27").Value = '  -2.23%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '37.24'
$ws.Range("E29").Value = '  -3.31%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '707.93'
$ws.Range("E30").Value = '  -1.45%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '13.46'
$ws.Range("E31").Value = '  -2.30%  '
$ws.Range("E32").Value = '  -1.16%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.90'
$ws.Range("E33").Value = '  -0.43%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '64.97'
$ws.Range("E34").Value = '  +5.84%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.451'
$ws.Range("E35").Value = '  +13.36%  '
$ws.Range("E36").Value = '  +0.60%  '
$ws.Range("E37").Value = '  -2.21%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '41.02'
$ws.Range("E38").Value = '  -2.99%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.151'
$ws.Range("E39").Value = '  +0.93%  '
$ws.Range("E40").Value = '  -0.13%  '
$ws.Range("E41").Value = '  -0.02%  '
$ws.Range("E42").Value = '  +1.77%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.92'
$ws.Range("E43").Value = '  -2.73%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.10'
$ws.Range("E44").Value = '  +5.36%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.07'
$ws.Range("E45").Value = '  -5.14%  '
$ws.Range("E46").Value = '  +1.48%  '
$ws.Range("E47").Value = '  +3.57%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.01'
$ws.Range("E48").Value = '  +6.05%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.41'
$ws.Range("E49").Value = '  -0.10%  '
$ws.Range("D50").Value = '0.0₆0352'
$ws.Range("E50").Value = '  -1.39%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.12'
$ws.Range("E51").Value = '  -1.60%  '
